$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cWidth = $ws.Range("C1").EntireColumn.ColumnWidth
Write-Host "C width before insert:" $cWidth
$ws.Range("D1:E1").EntireColumn.Insert()
$ws.Range("D1:E1").EntireColumn.ColumnWidth = $cWidth
$ws.Range("D1:E1").EntireColumn.Hidden = $true
Write-Host "Done"
